$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "#07-파이썬(Python) 반복문"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-07"

$ws.Range("D5").Value = "경계값 문제"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/28/Boundary_Value_Problem.html"

$ws.Range("D9").Value = "데이터 사이언스 마지막 강의 후기"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/data-science-class-review-finale/#utm_source=rss&utm_medium=rss&utm_campaign=data-science-class-review-finale"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 06월,『개인 맞춤형 미생물 유전체 분석』 교육생 모집 공고"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/400"

$ws.Range("D51").Value = "[MariaDB] Order by, 여러 개로 정렬하기 (더 중요한 것을 앞에)"
$ws.Range("E51").Value = "https://bskyvision.com/1204"
